$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'303.33"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'4.81%"
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(2, 7).Value = "'6"
$ws.Cells.Item(2, 7).Style = "Normal"
$ws.Cells.Item(3, 4).Value = "'34.81"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'12.32%"
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(3, 7).Value = "'6"
$ws.Cells.Item(3, 7).Style = "Normal"
$ws.Cells.Item(4, 4).Value = "'5.137"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "'4.34%"
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(4, 7).Value = "'6"
$ws.Cells.Item(4, 7).Style = "Normal"
$ws.Cells.Item(5, 4).Value = "'0.07743"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'4.68%"
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(5, 7).Value = "'6"
$ws.Cells.Item(5, 7).Style = "Normal"
$ws.Cells.Item(6, 4).Value = "'2.343"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'4.23%"
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(6, 7).Value = "'6"
$ws.Cells.Item(6, 7).Style = "Normal"
$ws.Cells.Item(7, 4).Value = "'8.016"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "'3.92%"
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(7, 7).Value = "'6"
$ws.Cells.Item(7, 7).Style = "Normal"
$ws.Cells.Item(8, 4).Value = "'3.952"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "'5.39%"
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(8, 7).Value = "'6"
$ws.Cells.Item(8, 7).Style = "Normal"
$ws.Cells.Item(9, 4).Value = "'0.9290"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "'1.95%"
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(9, 7).Value = "'6"
$ws.Cells.Item(9, 7).Style = "Normal"
$ws.Cells.Item(10, 4).Value = "'0.09952"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'12.83%"
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(10, 7).Value = "'6"
$ws.Cells.Item(10, 7).Style = "Normal"
$ws.Cells.Item(11, 4).Value = "'0.1798"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "'6.67%"
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(11, 7).Value = "'6"
$ws.Cells.Item(11, 7).Style = "Normal"
$ws.Cells.Item(12, 4).Value = "'0.08614"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "'4.54%"
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(12, 7).Value = "'6"
$ws.Cells.Item(12, 7).Style = "Normal"
$ws.Cells.Item(13, 4).Value = "'0.03318"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'6.29%"
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(13, 7).Value = "'6"
$ws.Cells.Item(13, 7).Style = "Normal"
$ws.Cells.Item(14, 4).Value = "'0.09895"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'-0.66%"
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(14, 7).Value = "'6"
$ws.Cells.Item(14, 7).Style = "Normal"
$ws.Cells.Item(15, 4).Value = "'0.001506"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'0.47%"
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(15, 7).Value = "'6"
$ws.Cells.Item(15, 7).Style = "Normal"
$ws.Cells.Item(16, 4).Value = "'0.005752"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'-1.51%"
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(16, 7).Value = "'6"
$ws.Cells.Item(16, 7).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'-0.82%"
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(17, 7).Value = "'6"
$ws.Cells.Item(17, 7).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'2.28%"
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(18, 7).Value = "'6"
$ws.Cells.Item(18, 7).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'1.17%"
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(19, 7).Value = "'6"
$ws.Cells.Item(19, 7).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'2.85%"
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(20, 7).Value = "'6"
$ws.Cells.Item(20, 7).Style = "Normal"
$ws.Cells.Item(21, 4).Value = "'4.346"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'12.49%"
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(21, 7).Value = "'6"
$ws.Cells.Item(21, 7).Style = "Normal"
$ws.Cells.Item(22, 4).Value = "'0.2389"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'9.15%"
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(22, 7).Value = "'6"
$ws.Cells.Item(22, 7).Style = "Normal"
$ws.Cells.Item(23, 4).Value = "'0.04568"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'0.34%"
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(23, 7).Value = "'6"
$ws.Cells.Item(23, 7).Style = "Normal"
$ws.Cells.Item(24, 4).Value = "'0.001219"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'0.73%"
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(24, 7).Value = "'6"
$ws.Cells.Item(24, 7).Style = "Normal"
$ws.Cells.Item(25, 4).Value = "'0.004458"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'-2.60%"
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(25, 7).Value = "'6"
$ws.Cells.Item(25, 7).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "'-0.07%"
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(26, 7).Value = "'6"
$ws.Cells.Item(26, 7).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'-0.20%"
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(27, 7).Value = "'6"
$ws.Cells.Item(27, 7).Style = "Normal"
$ws.Cells.Item(28, 7).Value = "'6"
$ws.Cells.Item(28, 7).Style = "Normal"
$ws.Cells.Item(29, 7).Value = "'6"
$ws.Cells.Item(29, 7).Style = "Normal"
$ws.Cells.Item(30, 7).Value = "'6"
$ws.Cells.Item(30, 7).Style = "Normal"
$ws.Cells.Item(31, 7).Value = "'6"
$ws.Cells.Item(31, 7).Style = "Normal"
$ws.Cells.Item(32, 7).Value = "'6"
$ws.Cells.Item(32, 7).Style = "Normal"
$ws.Cells.Item(33, 7).Value = "'6"
$ws.Cells.Item(33, 7).Style = "Normal"
$ws.Cells.Item(34, 7).Value = "'6"
$ws.Cells.Item(34, 7).Style = "Normal"
$ws.Cells.Item(35, 7).Value = "'6"
$ws.Cells.Item(35, 7).Style = "Normal"
$ws.Cells.Item(36, 7).Value = "'6"
$ws.Cells.Item(36, 7).Style = "Normal"
$ws.Cells.Item(37, 7).Value = "'6"
$ws.Cells.Item(37, 7).Style = "Normal"
$ws.Cells.Item(38, 7).Value = "'6"
$ws.Cells.Item(38, 7).Style = "Normal"
$ws.Cells.Item(39, 4).Value = "'0.01780"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'12.00%"
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(39, 7).Value = "'6"
$ws.Cells.Item(39, 7).Style = "Normal"
$ws.Cells.Item(40, 4).Value = "'0.04795"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'7.31%"
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(40, 7).Value = "'6"
$ws.Cells.Item(40, 7).Style = "Normal"
$ws.Cells.Item(41, 4).Value = "'0.007746"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'6.38%"
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(41, 7).Value = "'6"
$ws.Cells.Item(41, 7).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'6.37%"
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(42, 7).Value = "'6"
$ws.Cells.Item(42, 7).Style = "Normal"
$ws.Cells.Item(43, 4).Value = "'0.006860"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'-29.11%"
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(43, 7).Value = "'6"
$ws.Cells.Item(43, 7).Style = "Normal"
$ws.Cells.Item(44, 4).Value = "'0.002093"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'-6.27%"
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(44, 7).Value = "'6"
$ws.Cells.Item(44, 7).Style = "Normal"
$ws.Cells.Item(45, 4).Value = "'0.009181"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'3.42%"
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(45, 7).Value = "'6"
$ws.Cells.Item(45, 7).Style = "Normal"
$ws.Cells.Item(46, 4).Value = "'0.00006125"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'0.44%"
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(46, 7).Value = "'6"
$ws.Cells.Item(46, 7).Style = "Normal"
$ws.Cells.Item(47, 4).Value = "'0.00000000751"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'-0.07%"
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(47, 7).Value = "'6"
$ws.Cells.Item(47, 7).Style = "Normal"
$ws.Cells.Item(48, 4).Value = "'2.304"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'1.62%"
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(48, 7).Value = "'6"
$ws.Cells.Item(48, 7).Style = "Normal"
$ws.Cells.Item(49, 4).Value = "'0.002002"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'-0.07%"
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(49, 7).Value = "'6"
$ws.Cells.Item(49, 7).Style = "Normal"
$ws.Cells.Item(50, 4).Value = "'0.00002102"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'-0.07%"
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(50, 7).Value = "'6"
$ws.Cells.Item(50, 7).Style = "Normal"
$ws.Cells.Item(51, 4).Value = "'0.0002002"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "'-0.07%"
$ws.Cells.Item(51, 5).Style = "Normal"
$ws.Cells.Item(51, 7).Value = "'6"
$ws.Cells.Item(51, 7).Style = "Normal"
